# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Fri Oct 27 05:07:27 UTC 2023 with GitHub Actions".
# Only the cells that actually changed are touched; everything else is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "34.073.95"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "1.794.27"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'223.85"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'0.550"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'32.34"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "'0.285"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'0.0710"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "2.050.57"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "'10.98"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "1.779.13"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "34.098.16"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'4.17"
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "'68.05"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").Value = "'243.93"
$ws.Range("E19").Value = "  -3.95%  "
$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'10.75"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("E23").Value = "  -4.37%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("D25").Value = "'159.07"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").Value = "'16.29"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'3.67"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("E33").Value = "  -3.99%  "
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("D35").Value = "1.385.78"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "'0.650"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "'79.66"
$ws.Range("E39").Value = "  -6.44%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'2.35"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.917"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.71"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "0.0₆0139"
$ws.Range("E44").Value = "  +9.46%  "
$ws.Range("D45").Value = "'0.0500"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("D46").Value = "'107.58"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").Value = "1.952.04"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "'12.04"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  +0.00%  "
